$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting existing rows 52-60 down to 53-61
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with data (same record as the following
# row, but with updated date/price figures per the latest measurement)
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 44617
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100112021
$ws.Cells.Item(52, 7).Value = "Ají"
$ws.Cells.Item(52, 8).Value = "Inferno"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 120
$ws.Cells.Item(52, 11).Value = 13000
$ws.Cells.Item(52, 12).Value = 14000
$ws.Cells.Item(52, 13).Value = 13500
$ws.Cells.Item(52, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 900
$ws.Cells.Item(52, 17).Value = 15
$ws.Cells.Item(52, 18).Value = "Hortaliza"

Write-Output "done"
